$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header label text
$ws.Range("B1").Value = "Group"

# Round numeric stat values to fewer decimal places
$ws.Range("E2").Value = 261.528
$ws.Range("E3").Value = 544.577
$ws.Range("E4").Value = 14.8327
$ws.Range("E6").Value = 0.2864
